$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new blank row above row 18, pushing the remaining rows (old 18-36) down to 19-37.
$ws.Rows.Item(18).Insert()

# Fill in the data for the newly inserted row 18.
$ws.Range("A18").Value = "1.1/2.0"
$ws.Range("B18").Value = "Yale"
$ws.Range("D18").Value = 5368
$ws.Range("E18").Value = "low"
$ws.Range("F18").Value = 5
$ws.Range("H18").Value = "Implement a CAS plugin for the Authentication Service."

# Match the new row height recorded for row 18.
$ws.Rows.Item(18).RowHeight = 16.5

# Update the active selection to match the post-edit state.
$ws.Range("J8").Select()
